$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 13 (shifts old rows 13..23 down to 14..24).
#    All row heights and most label/value pairs shift down by one row and
#    keep their original formatting/content, which already matches the
#    target layout except for a handful of cells fixed below.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# Copy the value-column formatting (styles only) from the row below into the
# freshly inserted row 13 so that B13/C13 end up with the normal styles
# (s="2" / s="3") instead of inheriting the label style from row 12.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)

# Populate the new row 13 (professor's name, no label in column A).
$ws.Range("B13").Value = "6007846 - Júlio César dos Santos"
$ws.Range("C13").Value = "6007846 - Júlio César dos Santos"
$ws.Range("A13").Clear()

# ---------------------------------------------------------------------------
# 2. Fix the "Objetivos:" value row (row 10) - it used to hold the
#    professor's name by mistake; it must hold the actual objectives text.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "Desenvolver o aprendizado teórico e prático da Bioquímica através da execução de práticas de laboratório baseadas na evolução do conteúdo teórico ministrado na disciplina Bioquímica II."
$ws.Range("C10").Value = "Desenvolver o aprendizado teórico e prático da Bioquímica através da execução de práticas de laboratório baseadas na evolução do conteúdo teórico ministrado na disciplina Bioquímica II."

# ---------------------------------------------------------------------------
# 3. Row 14 ("Programa resumido:") - replace the bogus "Semestral" value
#    with the actual short syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "Propriedades gerais de glicídios; Fermentação anaeróbia; Extração deClorofila e Reação de Hill; Transporte de glicídios e indução de enzimas."
$ws.Range("C14").Value = "Propriedades gerais de glicídios; Fermentação anaeróbia; Extração deClorofila e Reação de Hill; Transporte de glicídios e indução de enzimas."

# ---------------------------------------------------------------------------
# 4. Row 16 ("Programa:") - replace the bogus date value with the actual
#    (long) syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "Propriedades gerais de glicídios:principais testes qualitativos para identificação e diferenciação de glicídios; aplicação de certas reações coloridas e dosagem espectrofotométrica de monossacarídeos redutores.Fermentação anaeróbia: Conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono;cálculo da eficiência do processo; ação de um inibidor da glicólise. Extração de clorofila e reação deHill: estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura e fase luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. Transporte de glicídios e indução de enzimas: conceitos gerais; enzimas do catabolismo da galactose; repressão, inativação emodificação catabólicas; sistemas enzimáticos constitutivos e induzidos em células de levedura"
$ws.Range("C16").Value = "Propriedades gerais de glicídios:principais testes qualitativos para identificação e diferenciação de glicídios; aplicação de certas reações coloridas e dosagem espectrofotométrica de monossacarídeos redutores.Fermentação anaeróbia: Conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono;cálculo da eficiência do processo; ação de um inibidor da glicólise. Extração de clorofila e reação deHill: estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura e fase luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. Transporte de glicídios e indução de enzimas: conceitos gerais; enzimas do catabolismo da galactose; repressão, inativação emodificação catabólicas; sistemas enzimáticos constitutivos e induzidos em células de levedura"

# ---------------------------------------------------------------------------
# 5. Row 19 ("Método:") - replace the bogus professor's name with the real
#    evaluation method text.
# ---------------------------------------------------------------------------
$ws.Range("B19").Value = "A avaliação será feita por meio de uma prova escrita e notas de relatórios (R)."
$ws.Range("C19").Value = "A avaliação será feita por meio de uma prova escrita e notas de relatórios (R)."

# ---------------------------------------------------------------------------
# 6. Row 20 ("Critério:") - shift the evaluation-method text out, put the
#    final-grade formula text in.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1*2 + R)/3."
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1*2 + R)/3."

# ---------------------------------------------------------------------------
# 7. Row 21 ("Norma de recuperação:") - put the recovery-exam rule text in.
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada pela fórmula: MR = (NF + PR)/2"

# ---------------------------------------------------------------------------
# 8. Row 22 ("Bibliografia:") - put the actual bibliography text in.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "CISTERNAS, J. R. Fundamentos de bioquímica experimental. São Paulo : Atheneu, 2005. ISBN: 9788573791075.NELSON, D. L., COX. M. M. Princípios de bioquímica de Lehninger. Porto Alegre : Artmed, 2011. ISBN: 9788536324180.VOET, D., VOET, J. G. Bioquímica. Porto Alegre : Artmed, 2013. ISBN: 9788582710043."
$ws.Range("C22").Value = "CISTERNAS, J. R. Fundamentos de bioquímica experimental. São Paulo : Atheneu, 2005. ISBN: 9788573791075.NELSON, D. L., COX. M. M. Princípios de bioquímica de Lehninger. Porto Alegre : Artmed, 2011. ISBN: 9788536324180.VOET, D., VOET, J. G. Bioquímica. Porto Alegre : Artmed, 2013. ISBN: 9788582710043."

# ---------------------------------------------------------------------------
# 9. NOTE on column widths: the source XML narrows the first <col> element
#    from "min=1 max=2" down to "min=1 max=1" (column A's width definition
#    used to redundantly span into column B). This is a pure no-op cleanup:
#    column B's effective width is already governed by its own, later,
#    more specific <col min="2" max="2" width="60.7109375".../> entry, so
#    every reader (Excel, openpyxl, ...) already resolves A=30.7109375,
#    B=60.7109375, C=60.7109375 both before and after that XML tweak.
#    The COM ColumnWidth setter in this runtime only accepts pixel-rounded
#    values (e.g. 31.5) and cannot reproduce the exact 30.7109375 figure,
#    and it always leaves a redundant extra <col> entry behind when a
#    multi-column range is narrowed. Since touching it would only trade an
#    invisible structural difference for a visible numeric one, the column
#    widths are intentionally left as-is here.
# ---------------------------------------------------------------------------
